$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "62.142.70"
Set-TextCell "E2" "  -0.90%  "
Set-TextCell "D3" "3.010.72"
Set-TextCell "E3" "  -0.26%  "
Set-TextCell "E4" "  +0.02%  "
Set-TextCell "D5" "593.64"
Set-TextCell "E5" "  +1.15%  "
Set-TextCell "D6" "147.07"
Set-TextCell "E6" "  -0.45%  "
Set-TextCell "E7" "  +0.04%  "
Set-TextCell "D8" "3.011.02"
Set-TextCell "E8" "  -0.22%  "
Set-TextCell "D9" "0.518"
Set-TextCell "E9" "  -1.94%  "
Set-TextCell "D10" "6.32"
Set-TextCell "E10" "  +7.81%  "
Set-TextCell "E11" "  -0.54%  "
Set-TextCell "D12" "0.457"
Set-TextCell "E12" "  -1.00%  "
Set-TextCell "D13" "0.0000232"
Set-TextCell "E13" "  +0.22%  "
Set-TextCell "D14" "34.37"
Set-TextCell "E14" "  -1.48%  "
Set-TextCell "E15" "  +2.57%  "
Set-TextCell "D16" "3.592.65"
Set-TextCell "E16" "  +2.34%  "
Set-TextCell "D17" "62.125.37"
Set-TextCell "E17" "  -0.77%  "
Set-TextCell "D18" "6.98"
Set-TextCell "E18" "  -1.95%  "
Set-TextCell "D19" "3.006.09"
Set-TextCell "E19" "  -0.35%  "
Set-TextCell "D20" "446.15"
Set-TextCell "E20" "  -3.08%  "
Set-TextCell "D21" "14.17"
Set-TextCell "E21" "  +0.97%  "
Set-TextCell "D22" "0.687"
Set-TextCell "E22" "  -0.70%  "
Set-TextCell "D23" "7.40"
Set-TextCell "E23" "  -0.70%  "
Set-TextCell "D24" "82.25"
Set-TextCell "E24" "  +0.59%  "
Set-TextCell "D25" "10.99"
Set-TextCell "E25" "  +9.33%  "
Set-TextCell "D26" "2.24"
Set-TextCell "E26" "  +0.57%  "
Set-TextCell "D27" "12.07"
Set-TextCell "E27" "  -2.28%  "
Set-TextCell "E28" "  -0.05%  "
Set-TextCell "E29" "  +1.69%  "
Set-TextCell "D31" "7.20"
Set-TextCell "E31" "  +1.92%  "
Set-TextCell "D32" "2.10"
Set-TextCell "E32" "  -0.72%  "
Set-TextCell "D33" "27.44"
Set-TextCell "E33" "  -2.56%  "
Set-TextCell "E34" "  +0.51%  "
Set-TextCell "D35" "0.0₃0849"
Set-TextCell "E35" "  +3.14%  "
Set-TextCell "D36" "1.03"
Set-TextCell "E36" "  -0.24%  "
Set-TextCell "E37" "  +0.50%  "
Set-TextCell "D38" "50.25"
Set-TextCell "E38" "  -0.37%  "
Set-TextCell "E39" "  -4.23%  "
Set-TextCell "D40" "9.02"
Set-TextCell "E40" "  -1.83%  "
Set-TextCell "E41" "  +1.05%  "
Set-TextCell "D42" "0.124"
Set-TextCell "E42" "  -0.03%  "
Set-TextCell "D43" "41.58"
Set-TextCell "E43" "  +11.40%  "
Set-TextCell "D44" "0.283"
Set-TextCell "E44" "  +4.33%  "
Set-TextCell "D45" "394.97"
Set-TextCell "E45" "  +0.31%  "
Set-TextCell "E46" "  -2.27%  "
Set-TextCell "D47" "2.724.88"
Set-TextCell "E47" "  -0.50%  "
Set-TextCell "D48" "134.11"
Set-TextCell "E48" "  +3.86%  "
Set-TextCell "D50" "2.18"
Set-TextCell "E50" "  -1.21%  "
Set-TextCell "E51" "  -1.73%  "
